# Auto-generated Excel COM-interop script
# Applies updated market-price snapshot values to the Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2262.6667
$ws.Range("J32").Value = 1894
$ws.Range("L32").Value = 1894
$ws.Range("N32").Value = -2546

$ws.Range("H88").Value = 2661.8
$ws.Range("I88").Value = 3180.8
$ws.Range("J88").Value = 2488.8
$ws.Range("K88").Value = 3180.8
$ws.Range("L88").Value = 2488.8
$ws.Range("M88").Value = -2774.8
$ws.Range("N88").Value = -3300.8

$ws.Range("H91").Value = 2661.8
$ws.Range("I91").Value = 3180.8
$ws.Range("J91").Value = 2488.8
$ws.Range("K91").Value = 3180.8
$ws.Range("L91").Value = 2488.8
$ws.Range("M91").Value = -1776.8
$ws.Range("N91").Value = -5296.8

$ws.Range("H110").Value = 100000
$ws.Range("J110").Value = 100000
$ws.Range("L110").Value = 100000
$ws.Range("N110").Value = -108180

$ws.Range("H117").Value = 75000
$ws.Range("J117").Value = 75000
$ws.Range("L117").Value = 75000
$ws.Range("N117").Value = -84178

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 2500.8909
$ws.Range("I132").Value = 1213.5714
$ws.Range("J132").Value = 6659.923
$ws.Range("K132").Value = 3640.7142
$ws.Range("L132").Value = 19979.769
$ws.Range("M132").Value = -1110.7142
$ws.Range("N132").Value = -25039.769

$ws.Range("H141").Value = 5668.1177
$ws.Range("I141").Value = 5628.625
$ws.Range("K141").Value = 16885.875
$ws.Range("M141").Value = -11705.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 88613.52
$ws.Range("I74").Value = 88613.52
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 88613.52
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -87739.52
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 88613.52
$ws.Range("I77").Value = 88613.52
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 443067.6
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -438699.6
$ws.Range("N77").ClearContents()

$ws.Range("H132").Value = 5400.1934
$ws.Range("I132").Value = 3889.0386
$ws.Range("J132").Value = 13258.2
$ws.Range("K132").Value = 11667.1158
$ws.Range("L132").Value = 39774.60000000001
$ws.Range("M132").Value = -9137.1158
$ws.Range("N132").Value = -44834.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 24999
$ws.Range("J92").Value = 24999
$ws.Range("L92").Value = 24999
$ws.Range("N92").Value = -29991

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 250006370
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 250006370
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 250006370
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -250006594

$ws.Range("H7").Value = 157.33333
$ws.Range("I7").Value = 243
$ws.Range("K7").Value = 243
$ws.Range("M7").Value = -130

$ws.Range("H23").Value = 9996.666999999999
$ws.Range("I23").Value = 20000
$ws.Range("J23").Value = 4995
$ws.Range("K23").Value = 20000
$ws.Range("L23").Value = 4995
$ws.Range("M23").Value = -19760
$ws.Range("N23").Value = -5475

$ws.Range("H27").Value = 9996.666999999999
$ws.Range("I27").Value = 20000
$ws.Range("J27").Value = 4995
$ws.Range("K27").Value = 20000
$ws.Range("L27").Value = 4995
$ws.Range("M27").Value = -19808
$ws.Range("N27").Value = -5379

$ws.Range("H31").Value = 143845.34
$ws.Range("I31").Value = 183832.81
$ws.Range("K31").Value = 183832.81
$ws.Range("M31").Value = -183537.81

$ws.Range("H34").Value = 143845.34
$ws.Range("I34").Value = 183832.81
$ws.Range("K34").Value = 183832.81
$ws.Range("M34").Value = -183630.81

$ws.Range("H95").Value = 8925.4
$ws.Range("J95").Value = 8925.4
$ws.Range("L95").Value = 8925.4
$ws.Range("N95").Value = -14417.4

$ws.Range("H132").Value = 3132.5386
$ws.Range("I132").Value = 1975.7646
$ws.Range("J132").Value = 10998.6
$ws.Range("K132").Value = 5927.293799999999
$ws.Range("L132").Value = 32995.8
$ws.Range("M132").Value = -3397.293799999999
$ws.Range("N132").Value = -38055.8

$ws.Range("H134").Value = 3806.0425
$ws.Range("I134").Value = 3884.025
$ws.Range("J134").Value = 3360.4285
$ws.Range("K134").Value = 11652.075
$ws.Range("L134").Value = 10081.2855
$ws.Range("M134").Value = -9117.075000000001
$ws.Range("N134").Value = -15151.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1200
$ws.Range("J17").Value = 1400
$ws.Range("L17").Value = 4200
$ws.Range("N17").Value = -4538

$ws.Range("H39").Value = 1698.1428
$ws.Range("I39").Value = 815.7778
$ws.Range("J39").Value = 2359.9167
$ws.Range("K39").Value = 2447.3334
$ws.Range("L39").Value = 7079.750100000001
$ws.Range("M39").Value = -2153.3334
$ws.Range("N39").Value = -7667.750100000001

$ws.Range("H99").Value = 5114.875
$ws.Range("I99").Value = 4674.143
$ws.Range("K99").Value = 14022.429
$ws.Range("M99").Value = -11776.429

$ws.Range("H118").Value = 324.5
$ws.Range("I118").Value = 324.5
$ws.Range("K118").Value = 973.5
$ws.Range("M118").Value = 269.5

$ws.Range("H121").Value = 2354.3572
$ws.Range("I121").Value = 1688.1666
$ws.Range("J121").Value = 2854
$ws.Range("K121").Value = 5064.4998
$ws.Range("L121").Value = 8562
$ws.Range("M121").Value = -3754.4998
$ws.Range("N121").Value = -11182

$ws.Range("H134").Value = 4851.727
$ws.Range("I134").Value = 5358.625
$ws.Range("K134").Value = 16075.875
$ws.Range("M134").Value = -11005.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10674.272
$ws.Range("I70").Value = 10674.272
$ws.Range("K70").Value = 10674.272
$ws.Range("M70").Value = -10404.272

$ws.Range("H73").Value = 10674.272
$ws.Range("I73").Value = 10674.272
$ws.Range("K73").Value = 10674.272
$ws.Range("M73").Value = -9738.272000000001

$ws.Range("H80").Value = 4461.727
$ws.Range("I80").Value = 4119.8887
$ws.Range("J80").Value = 6000
$ws.Range("K80").Value = 4119.8887
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = -3121.8887
$ws.Range("N80").Value = -7996

$ws.Range("H83").Value = 4461.727
$ws.Range("I83").Value = 4119.8887
$ws.Range("J83").Value = 6000
$ws.Range("K83").Value = 20599.4435
$ws.Range("L83").Value = 30000
$ws.Range("M83").Value = -15607.4435
$ws.Range("N83").Value = -39984

$ws.Range("H122").Value = 2228.889
$ws.Range("I122").Value = 2071.8235
$ws.Range("J122").Value = 4899
$ws.Range("K122").Value = 6215.470499999999
$ws.Range("L122").Value = 14697
$ws.Range("M122").Value = -3765.470499999999
$ws.Range("N122").Value = -19597

$ws.Range("H132").Value = 49669.695
$ws.Range("I132").Value = 62735.707
$ws.Range("K132").Value = 188207.121
$ws.Range("M132").Value = -185677.121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 18923.2
$ws.Range("J105").Value = 18923.2
$ws.Range("L105").Value = 18923.2
$ws.Range("N105").Value = -25911.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 616.95
$ws.Range("I107").Value = 583.5333000000001
$ws.Range("J107").Value = 717.2
$ws.Range("K107").Value = 1750.5999
$ws.Range("L107").Value = 2151.6
$ws.Range("M107").Value = 169.4000999999998
$ws.Range("N107").Value = -5991.6

$ws.Range("H124").Value = 70429
$ws.Range("J124").Value = 70429
$ws.Range("L124").Value = 70429
$ws.Range("N124").Value = -80249

$ws.Range("H135").Value = 67048.125
$ws.Range("J135").Value = 67048.125
$ws.Range("L135").Value = 67048.125
$ws.Range("N135").Value = -77188.125

